$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" worksheet right before the "总计" (Total) sheet.
#    Seed it with the same look (header style, empty A1) as the other
#    quarterly sheets by copying the formatting from "2021-Q1" first, then
#    overwrite the cell values/content for the new quarter.
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q1"

# NB: the host aliases the "Before" sheet argument onto the freshly-created
# sheet, so $beforeSheet (and the original $totalSheet handle) now points at
# the NEW "2022-Q1" tab, not at "总计" any more. Re-resolve "总计" by name
# after the Add() call to get a handle on the real sheet.
$totalSheet = $wb.Worksheets.Item("总计")

$srcSheet = $wb.Worksheets.Item("2021-Q1")
$srcSheet.Range("A1:H3").Copy($newSheet.Range("A1:H3"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "004194"

$newSheet.Range("C2").Value = "招商中证1000指数增强A"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.76"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "94.40"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "1.09"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0192"

$newSheet.Range("H2").Value = 4

# Row 3
$newSheet.Range("A3").Value = 1

$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B3").Value = "004195"

$newSheet.Range("C3").Value = "招商中证1000指数增强C"

$newSheet.Range("D3").NumberFormat = "@"
$newSheet.Range("D3").Value = "0.68"

$newSheet.Range("E3").NumberFormat = "@"
$newSheet.Range("E3").Value = "94.40"

$newSheet.Range("F3").NumberFormat = "@"
$newSheet.Range("F3").Value = "1.09"

$newSheet.Range("G3").NumberFormat = "@"
$newSheet.Range("G3").Value = "0.0074"

$newSheet.Range("H3").Value = 4

# ---------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" (Total) summary sheet, pushing the
#    existing quarters down by one row and renumbering the running index in
#    column A.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.03

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

Write-Output "done"
